$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 664.11536
$ws.Range("I28").Value = 512.0952
$ws.Range("K28").Value = 512.0952
$ws.Range("M28").Value = -27.09519999999998
$ws.Range("H69").Value = 690999.75
$ws.Range("J69").Value = 690999.75
$ws.Range("L69").Value = 2072999.25
$ws.Range("N69").Value = -2074747.25
$ws.Range("H72").Value = 690999.75
$ws.Range("J72").Value = 690999.75
$ws.Range("L72").Value = 6218997.75
$ws.Range("N72").Value = -6227733.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4623.9697
$ws.Range("I61").Value = 748.6667
$ws.Range("K61").Value = 748.6667
$ws.Range("M61").Value = -536.6667
$ws.Range("H97").Value = 1479.5312
$ws.Range("I97").Value = 1175
$ws.Range("K97").Value = 1175
$ws.Range("M97").Value = -679
$ws.Range("H136").Value = 4623.9697
$ws.Range("I136").Value = 748.6667
$ws.Range("K136").Value = 2246.0001
$ws.Range("M136").Value = 303.9998999999998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1356.5
$ws.Range("I86").Value = 1311.4482
$ws.Range("K86").Value = 1311.4482
$ws.Range("M86").Value = -188.4482
$ws.Range("H89").Value = 1356.5
$ws.Range("I89").Value = 1311.4482
$ws.Range("K89").Value = 6557.241
$ws.Range("M89").Value = -941.241
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 8351
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H10").Value = 452.5
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 141632290
$ws.Range("J12").Value = 30267000
$ws.Range("L12").Value = 30267000
$ws.Range("N12").Value = -30267340
$ws.Range("H14").Value = 17297.25
$ws.Range("J14").Value = 17297.25
$ws.Range("L14").Value = 17297.25
$ws.Range("N14").Value = -17637.25
$ws.Range("H15").Value = 533
$ws.Range("J15").Value = 560.6667
$ws.Range("L15").Value = 560.6667
$ws.Range("N15").Value = -900.6667
$ws.Range("H26").Value = 4950
$ws.Range("J26").Value = 4950
$ws.Range("L26").Value = 4950
$ws.Range("N26").Value = -5524
$ws.Range("H29").Value = 23500
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 23500
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 23500
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -24086
$ws.Range("H86").Value = 13805.046
$ws.Range("I86").Value = 16490.75
$ws.Range("J86").Value = 10582.2
$ws.Range("K86").Value = 16490.75
$ws.Range("L86").Value = 10582.2
$ws.Range("M86").Value = -15367.75
$ws.Range("N86").Value = -12828.2
$ws.Range("H89").Value = 13805.046
$ws.Range("I89").Value = 16490.75
$ws.Range("J89").Value = 10582.2
$ws.Range("K89").Value = 82453.75
$ws.Range("L89").Value = 52911
$ws.Range("M89").Value = -76837.75
$ws.Range("N89").Value = -64143
$ws.Range("H94").Value = 1543.5416
$ws.Range("I94").Value = 1402.75
$ws.Range("J94").Value = 1613.9375
$ws.Range("K94").Value = 1402.75
$ws.Range("L94").Value = 1613.9375
$ws.Range("M94").Value = -951.75
$ws.Range("N94").Value = -2515.9375
$ws.Range("H134").Value = 2209.8
$ws.Range("I134").Value = 1839.5454
$ws.Range("K134").Value = 5518.6362
$ws.Range("M134").Value = -2983.6362
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 614.6667
$ws.Range("I36").Value = 614.6667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1844.0001
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1675.0001
$ws.Range("N36").ClearContents()
$ws.Range("H76").Value = 4578.9473
$ws.Range("J76").Value = 4805.5557
$ws.Range("L76").Value = 14416.6671
$ws.Range("N76").Value = -15182.6671
$ws.Range("H79").Value = 4578.9473
$ws.Range("J79").Value = 4805.5557
$ws.Range("L79").Value = 14416.6671
$ws.Range("N79").Value = -17068.6671
$ws.Range("H109").Value = 1307.1
$ws.Range("I109").Value = 179.16667
$ws.Range("K109").Value = 537.50001
$ws.Range("M109").Value = 502.49999
$ws.Range("H113").Value = 733.3333
$ws.Range("J113").Value = 733.3333
$ws.Range("L113").Value = 2199.9999
$ws.Range("N113").Value = -6539.9999
$ws.Range("H132").Value = 2313.4285
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3412.5293
$ws.Range("I7").Value = 3696.111
$ws.Range("J7").Value = 3093.5
$ws.Range("K7").Value = 3696.111
$ws.Range("L7").Value = 3093.5
$ws.Range("M7").Value = -3584.111
$ws.Range("N7").Value = -3317.5
$ws.Range("H22").Value = 1191.0769
$ws.Range("I22").Value = 1039.2667
$ws.Range("K22").Value = 1039.2667
$ws.Range("M22").Value = -744.2666999999999
$ws.Range("H27").Value = 1191.0769
$ws.Range("I27").Value = 1039.2667
$ws.Range("K27").Value = 1039.2667
$ws.Range("M27").Value = -932.2666999999999
$ws.Range("H46").Value = 2674.889
$ws.Range("I46").Value = 786.2143
$ws.Range("J46").Value = 4708.846
$ws.Range("K46").Value = 786.2143
$ws.Range("L46").Value = 4708.846
$ws.Range("M46").Value = -598.2143
$ws.Range("N46").Value = -5084.846
$ws.Range("H55").Value = 1076.1428
$ws.Range("I55").Value = 265.625
$ws.Range("K55").Value = 265.625
$ws.Range("M55").Value = -92.625
$ws.Range("H100").Value = 3333.5
$ws.Range("I100").Value = 3000.2
$ws.Range("K100").Value = 3000.2
$ws.Range("M100").Value = -2459.2
$ws.Range("H126").Value = 3412.5293
$ws.Range("I126").Value = 3696.111
$ws.Range("J126").Value = 3093.5
$ws.Range("K126").Value = 11088.333
$ws.Range("L126").Value = 9280.5
$ws.Range("M126").Value = -8618.332999999999
$ws.Range("N126").Value = -14220.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 74829
$ws.Range("I122").Value = 82807.625
$ws.Range("K122").Value = 248422.875
$ws.Range("M122").Value = -245972.875
$ws.Range("H136").Value = 26733.666
$ws.Range("I136").Value = 29047.79
$ws.Range("K136").Value = 87143.37
$ws.Range("M136").Value = -84593.37
